$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.664.61"
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").Value = "1.598.81"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("E6").Value = "  +1.02%  "
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +0.64%  "
$ws.Range("D10").Value = "'19.59"
$ws.Range("E10").Value = "  -0.90%  "
$ws.Range("D11").Value = "'0.0837"
$ws.Range("E11").Value = "  +0.03%  "
$ws.Range("D12").Value = "1.823.42"
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.03"
$ws.Range("E13").Value = "  +0.02%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.562.26"
$ws.Range("E14").Value = "  -2.57%  "
$ws.Range("E15").Value = "  +0.12%  "
$ws.Range("D16").Value = "'65.21"
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("D17").Value = "26.660.42"
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("E18").Value = "  +1.40%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "'209.46"
$ws.Range("E19").Value = "  -0.17%  "
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").Value = "'1.00"
$ws.Range("E20").Value = "  +0.30%  "
$ws.Range("D21").Value = "'7.06"
$ws.Range("E21").Value = "  +4.56%  "
$ws.Range("E22").Value = "  +0.62%  "
$ws.Range("E23").Value = "  +0.64%  "
$ws.Range("D24").Value = "'8.99"
$ws.Range("E24").Value = "  +0.78%  "
$ws.Range("D25").Value = "'145.32"
$ws.Range("E25").Value = "  -1.01%  "
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("E27").Value = "  -0.97%  "
$ws.Range("E28").Value = "  -0.43%  "
$ws.Range("D29").Value = "'15.30"
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("E30").Value = "  +2.06%  "
$ws.Range("E31").Value = "  +0.16%  "
$ws.Range("E32").Value = "  +0.27%  "
$ws.Range("E33").Value = "  +1.35%  "
$ws.Range("D34").Value = "1.282.67"
$ws.Range("E34").Value = "  -1.32%  "
$ws.Range("E35").Value = "  -7.76%  "
$ws.Range("E36").Value = "  +0.53%  "
$ws.Range("E37").Value = "  +1.02%  "
$ws.Range("E38").Value = "  -0.99%  "
$ws.Range("E39").Value = "  -1.05%  "
$ws.Range("E40").Value = "  +19.13%  "
$ws.Range("D41").Value = "'5.51"
$ws.Range("E41").Value = "  +2.56%  "
$ws.Range("D43").Value = "'0.784"
$ws.Range("E43").Value = "  -0.95%  "
$ws.Range("D44").Value = "'63.96"
$ws.Range("E44").Value = "  +0.20%  "
$ws.Range("D45").Value = "1.735.83"
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("D46").Value = "'91.51"
$ws.Range("E46").Value = "  +1.53%  "
$ws.Range("D47").Value = "'1.58"
$ws.Range("E47").Value = "  -3.06%  "
$ws.Range("E48").Value = "  +3.27%  "
$ws.Range("E49").Value = "  +0.55%  "
$ws.Range("E50").Value = "  -0.10%  "
$ws.Range("D51").Value = "'7.40"
$ws.Range("E51").Value = "  -1.57%  "
